$wb = $excel.ActiveWorkbook
$ws4 = $wb.Worksheets.Item("convertSeqFiles")

# Row 2 updates
$ws4.Range("D2").Value = $null
$ws4.Range("G2").Value = 400
$ws4.Range("H2").Value = 550

# Row 3: clear all data (A3:H3)
$ws4.Range("A3:H3").ClearContents()

# Row 4: clear all data (A4:H4)
$ws4.Range("A4:H4").ClearContents()

$ws6 = $wb.Worksheets.Item("deleteImageSequence")
$ws6.Range("A2").Value = "/scratch/utkur/utkarsh/HeatingHolder/Al-Film/05nm/20190718/20190718_16-12-48.635_converted/png"
$ws6.Range("C2").Value = 11899
$ws6.Range("A3:D6").Delete()
$ws6.Columns("A").WrapText = $true

$ws4.Activate()
$ws4.Range("G2").Select()
